$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Final Paper" column header (G1) — new shared string entry.
$ws.Range("G1").Value = "Final Paper"

# New "Final Paper" scores for the students who have one on record.
$ws.Range("G12").Formula = "=75/100"
$ws.Range("G14").Formula = "=90/100"
$ws.Range("G27").Formula = "=95/100"

# Size the new column roughly like the source workbook (~9.58 chars wide).
$ws.Columns("G").ColumnWidth = 8.65

# Move the active selection, like the editor left it.
[void]$ws.Range("G13").Select()
